# Add data for 2022-05-05 (update "through May 04" -> "through May 05")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its header label to reflect the new "through" date
$ws.Name = "Through 2022-05-05"
$ws.Range("B1").Value = "May 2022 (through May 05)"

# Helper to add an increment to a cell (treats blank as 0)
function Add-Value {
    param(
        [object]$Worksheet,
        [string]$CellRef,
        [double]$Increment
    )
    $cell = $Worksheet.Range($CellRef)
    $current = $cell.Value2
    if ($null -eq $current) {
        $current = 0
    }
    $cell.Value = $current + $Increment
}

# Englewood (row 2)
Add-Value $ws "B2" 1
Add-Value $ws "L2" 1

# Austin (row 3)
Add-Value $ws "AA3" 1

# Garfield Park (row 5)
Add-Value $ws "B5" 1
Add-Value $ws "AA5" 1
Add-Value $ws "AF5" 1

# South Shore (row 8)
Add-Value $ws "G8" 1
Add-Value $ws "AK8" 1

# Lake View (row 15)
Add-Value $ws "V15" 1

# Brighton Park (row 22)
Add-Value $ws "B22" 1

# Grand Crossing (row 23)
Add-Value $ws "G23" 1

# South Deering (row 34)
Add-Value $ws "G34" 1

# Avondale (row 35)
Add-Value $ws "G35" 1

# Douglas (row 38)
Add-Value $ws "G38" 1

# New City (row 39)
Add-Value $ws "L39" 1

# Morgan Park (row 41)
Add-Value $ws "G41" 1

# Jackson Park (row 70)
Add-Value $ws "G70" 1

# Washington Park (row 91)
Add-Value $ws "B91" 1

# West Lawn (row 93)
Add-Value $ws "B93" 1
